# Generate Report for Handback
# Update timestamp values recorded during handback report generation.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first file
# (this same timestamp string is also shown on the "de-de" sheet's
# Correspond Handoff Datetime for that same file, since it was identical).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 09:08:22"

# "zh-cn" sheet: handoff / handback datetimes for the first file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 09:08:16"
$wsZhCn.Range("K2").Value = "2016-08-18 09:08:45"

# "de-de" sheet: handback datetime for the first file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-18 09:08:53"
